# Add "total_data" (H) and "publish_date" (I) columns to the victims sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - copy the header style used by the existing G1 header cell
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("H1").Value2 = "total_data"
$ws.Range("I1").Value2 = "publish_date"

$data = @(
    @{ Row = 2;  Total = "70 GB";     Date = "25 Sep 2023" },
    @{ Row = 3;  Total = "80.68 KB";  Date = "25 Sep 2023" },
    @{ Row = 4;  Total = "191.62 KB"; Date = "20 Sep 2023" },
    @{ Row = 5;  Total = "1.64 MB";   Date = "17 Sep 2023" },
    @{ Row = 6;  Total = "804.51 KB"; Date = "16 Sep 2023" },
    @{ Row = 7;  Total = "594.16 KB"; Date = "14 Sep 2023" },
    @{ Row = 8;  Total = "496.55 KB"; Date = "10 Sep 2023" },
    @{ Row = 9;  Total = "415.42 KB"; Date = "10 Sep 2023" },
    @{ Row = 10; Total = "873.62 KB"; Date = "21 Sep 2023" },
    @{ Row = 11; Total = "17 GB";     Date = "24 Sep 2023" },
    @{ Row = 12; Total = "25 GB";     Date = "24 Sep 2023" },
    @{ Row = 13; Total = "20 GB";     Date = "24 Sep 2023" },
    @{ Row = 14; Total = "218.02 KB"; Date = "20 Sep 2023" }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 8).Value2 = $entry.Total
    $ws.Cells.Item($r, 9).Value2 = $entry.Date
}
